$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the SKU number list in column B (rows 2-13), keep their existing style.
$ws.Range("B2:B13").ClearContents()

# B14 held the stray shared string (" 469137733") - clear it too; its slot in
# the shared-string table gets reused below for the new URL text.
$ws.Range("B14").ClearContents()

# E2 becomes the product-card link.
$ws.Range("E2").Value = "https://www.ozon.ru/brand/karcher-26303230/"

# Highlight G2 with a yellow fill (fresh/clean style - clear first then fill).
$ws.Range("G2").Clear()
$ws.Range("G2").Interior.Color = 65535

# Selection moves to E2.
$ws.Range("E2").Select()
